# "last commit of the day"
# Rename the generic "Load Switch IC" line item (row 11) to be voltage-specific
# ("5V Load Switch IC"), and add two new BOM rows:
#   row 13 - a 9V load-switch IC (SI1865DDL-T1-BE3)
#   row 14 - a relay (G6DN-1A DC4.5) used as a 5A 4.5VDC relay
# Finally move the active selection to F7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: clarify the existing "Load Switch IC" part as the 5V variant.
$ws.Range("B11").Value = "5V Load Switch IC"

# Row 13: new 9V load switch IC line item.
$ws.Range("B13").Value = "9V Load Switch IC"
$ws.Range("F13").Value = "78-SI1865DDL-T1-BE3 "
$ws.Range("C13").Value = "SI1865DDL-T1-BE3 "
$ws.Range("E13").Value = "Mouser"

# Row 14: new relay line item.
$ws.Range("F14").Value = "653-G6DN1ADC45"
$ws.Range("C14").Value = "G6DN-1A DC4.5"
$ws.Range("B14").Value = "Relais 5A 4.5VDC"
$ws.Range("E14").Value = "Mouser"

# Move the selection as recorded in the saved view state.
$null = $ws.Range("F7").Select()
